$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-70 down to 51-71.
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the new "Haba" price record.
$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 44784
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = 100112026
$ws.Range("G50").Value = "Haba"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = 10000
$ws.Range("N50").Value = "$/saco 25 kilos"
$ws.Range("O50").Value = "Provincia del Elquí"
$ws.Range("P50").Value = 400
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
